$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref="D2"; Val="307.46"},
    @{Ref="E2"; Val="0.03%"},
    @{Ref="E3"; Val="0.14%"},
    @{Ref="D4"; Val="5.241"},
    @{Ref="E4"; Val="2.06%"},
    @{Ref="D5"; Val="0.07663"},
    @{Ref="E5"; Val="0.76%"},
    @{Ref="D6"; Val="1.630"},
    @{Ref="E6"; Val="0.09%"},
    @{Ref="D7"; Val="0.9166"},
    @{Ref="E7"; Val="1.79%"},
    @{Ref="D8"; Val="2.444"},
    @{Ref="E8"; Val="-0.81%"},
    @{Ref="D9"; Val="0.1253"},
    @{Ref="E9"; Val="12.10%"},
    @{Ref="D10"; Val="0.1836"},
    @{Ref="E10"; Val="3.71%"},
    @{Ref="D11"; Val="0.09250"},
    @{Ref="E11"; Val="0.18%"},
    @{Ref="D12"; Val="0.04279"},
    @{Ref="E12"; Val="2.18%"},
    @{Ref="D13"; Val="0.1051"},
    @{Ref="E13"; Val="0.27%"},
    @{Ref="D14"; Val="0.001261"},
    @{Ref="E14"; Val="0.75%"},
    @{Ref="D15"; Val="0.005835"},
    @{Ref="E15"; Val="0.48%"},
    @{Ref="D17"; Val="3.354"},
    @{Ref="E17"; Val="-0.08%"},
    @{Ref="D18"; Val="4.325"},
    @{Ref="E18"; Val="2.08%"},
    @{Ref="D19"; Val="0.3335"},
    @{Ref="D20"; Val="7.144"},
    @{Ref="E20"; Val="9.01%"},
    @{Ref="E21"; Val="1.55%"},
    @{Ref="E22"; Val="8.06%"},
    @{Ref="D23"; Val="0.04078"},
    @{Ref="E23"; Val="-1.44%"},
    @{Ref="D24"; Val="0.001262"},
    @{Ref="E24"; Val="3.20%"},
    @{Ref="D25"; Val="0.004150"},
    @{Ref="E25"; Val="3.82%"},
    @{Ref="E26"; Val="-2.02%"},
    @{Ref="D38"; Val="0.02471"},
    @{Ref="E38"; Val="3.14%"},
    @{Ref="D39"; Val="0.05279"},
    @{Ref="E39"; Val="1.69%"},
    @{Ref="D40"; Val="0.007851"},
    @{Ref="E40"; Val="0.91%"},
    @{Ref="D41"; Val="0.1315"},
    @{Ref="E41"; Val="1.18%"},
    @{Ref="E42"; Val="-1.86%"},
    @{Ref="D43"; Val="0.001917"},
    @{Ref="E43"; Val="-2.76%"},
    @{Ref="E44"; Val="2.58%"},
    @{Ref="D45"; Val="0.3051"},
    @{Ref="E45"; Val="-0.19%"},
    @{Ref="D46"; Val="0.00006748"},
    @{Ref="E46"; Val="0.41%"},
    @{Ref="D47"; Val="0.00000000752"},
    @{Ref="E47"; Val="0.27%"},
    @{Ref="D48"; Val="0.2055"},
    @{Ref="E48"; Val="2,198.84%"},
    @{Ref="E49"; Val="-2.42%"},
    @{Ref="E50"; Val="0.27%"},
    @{Ref="E51"; Val="0.27%"}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $u.Val
    $cell.Style = $originalStyle
}
